$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JPMGE")

# Update recomputed result values (rows 2-9, 17, 18)
$ws.Cells.Item(2, 5).Value = 6.0271570694089363
$ws.Cells.Item(2, 6).Value = 3.9819768629877865
$ws.Cells.Item(2, 7).Value = 6.0271570813322075
$ws.Cells.Item(2, 8).Value = 3.9819768428170277
$ws.Cells.Item(2, 9).Value = 6.0271570694027758
$ws.Cells.Item(2, 10).Value = 3.9819768629948684
$ws.Cells.Item(2, 11).Value = 6.0285029175131495
$ws.Cells.Item(2, 12).Value = 3.9810878966780647
$ws.Cells.Item(2, 13).Value = 6.0285029321430184
$ws.Cells.Item(2, 14).Value = 3.9810878724541241
$ws.Cells.Item(2, 15).Value = 6.0285029175126947
$ws.Cells.Item(2, 16).Value = 3.9810878966777277
$ws.Cells.Item(2, 17).Value = 6.0259370963253147
$ws.Cells.Item(2, 18).Value = 3.9827830288230945
$ws.Cells.Item(2, 19).Value = 6.025937105223008
$ws.Cells.Item(2, 20).Value = 3.9827830128948642
$ws.Cells.Item(2, 21).Value = 6.0259370963252268
$ws.Cells.Item(2, 22).Value = 3.9827830288229387
$ws.Cells.Item(3, 5).Value = 2.0306618580514879
$ws.Cells.Item(3, 6).Value = 1.9698011188502622
$ws.Cells.Item(3, 7).Value = 2.0306618717316747
$ws.Cells.Item(3, 8).Value = 1.9698011007128038
$ws.Cells.Item(3, 9).Value = 2.0306618580444193
$ws.Cells.Item(3, 10).Value = 1.9698011188587803
$ws.Cells.Item(3, 11).Value = 2.0322072644723499
$ws.Cells.Item(3, 12).Value = 1.968303169627029
$ws.Cells.Item(3, 13).Value = 2.0322072812850394
$ws.Cells.Item(3, 14).Value = 1.9683031476035679
$ws.Cells.Item(3, 15).Value = 2.0322072644718245
$ws.Cells.Item(3, 16).Value = 1.9683031696273214
$ws.Cells.Item(3, 17).Value = 2.0292631412455764
$ws.Cells.Item(3, 18).Value = 1.9711588500763715
$ws.Cells.Item(3, 19).Value = 2.0292631514394697
$ws.Cells.Item(3, 20).Value = 1.9711588362391685
$ws.Cells.Item(3, 21).Value = 2.0292631412454765
$ws.Cells.Item(3, 22).Value = 1.9711588500763977
$ws.Cells.Item(4, 5).Value = 1.9725936953390322
$ws.Cells.Item(4, 6).Value = 2.0277870751833285
$ws.Cells.Item(4, 7).Value = 1.9725936831953861
$ws.Cells.Item(4, 8).Value = 2.0277870814064576
$ws.Cells.Item(4, 9).Value = 1.9725936953453067
$ws.Cells.Item(4, 10).Value = 2.0277870751784119
$ws.Cells.Item(4, 11).Value = 1.9712223427385891
$ws.Cells.Item(4, 12).Value = 2.0291977791012972
$ws.Cells.Item(4, 13).Value = 1.9712223278246288
$ws.Cells.Item(4, 14).Value = 2.0291977870311428
$ws.Cells.Item(4, 15).Value = 1.9712223427390549
$ws.Cells.Item(4, 16).Value = 2.029197779100492
$ws.Cells.Item(4, 17).Value = 1.9738356999448601
$ws.Cells.Item(4, 18).Value = 2.026511122537574
$ws.Cells.Item(4, 19).Value = 1.9738356908902481
$ws.Cells.Item(4, 20).Value = 2.0265111267215343
$ws.Cells.Item(4, 21).Value = 1.9738356999449487
$ws.Cells.Item(4, 22).Value = 2.0265111225373742
$ws.Cells.Item(5, 5).Value = 9.9690552256106191
$ws.Cells.Item(5, 6).Value = 6.018624497716182
$ws.Cells.Item(5, 7).Value = 9.9690552116775866
$ws.Cells.Item(5, 8).Value = 6.0186244912562668
$ws.Cells.Item(5, 9).Value = 9.9690552256178169
$ws.Cells.Item(5, 10).Value = 6.018624497716913
$ws.Cells.Item(5, 11).Value = 9.9674805327690024
$ws.Cells.Item(5, 12).Value = 6.019575338295823
$ws.Cells.Item(5, 13).Value = 9.9674805156298323
$ws.Cells.Item(5, 14).Value = 6.0195753310938045
$ws.Cells.Item(5, 15).Value = 9.9674805327695388
$ws.Cells.Item(5, 16).Value = 6.0195753382948372
$ws.Cells.Item(5, 17).Value = 9.9704792020198632
$ws.Cells.Item(5, 18).Value = 6.0177649222561822
$ws.Cells.Item(5, 19).Value = 9.9704791916461968
$ws.Cells.Item(5, 20).Value = 6.0177649165035128
$ws.Cells.Item(5, 21).Value = 9.9704792020199662
$ws.Cells.Item(5, 22).Value = 6.0177649222559033
$ws.Cells.Item(6, 5).Value = 1.0326001216556313
$ws.Cells.Item(6, 7).Value = 1.0326001357507444
$ws.Cells.Item(6, 9).Value = 1.032600121646702
$ws.Cells.Item(6, 11).Value = 1.0342646526340822
$ws.Cells.Item(6, 13).Value = 1.0342646699325051
$ws.Cells.Item(6, 15).Value = 1.0342646526333701
$ws.Cells.Item(6, 17).Value = 1.0310953593284131
$ws.Cells.Item(6, 19).Value = 1.031095369634073
$ws.Cells.Item(6, 21).Value = 1.0310953593282555
$ws.Cells.Item(7, 5).Value = 3.0648347968089404
$ws.Cells.Item(7, 7).Value = 3.0648348259489921
$ws.Cells.Item(7, 9).Value = 3.0648347967902407
$ws.Cells.Item(7, 11).Value = 3.068125768905853
$ws.Cells.Item(7, 13).Value = 3.0681258045606343
$ws.Cells.Item(7, 15).Value = 3.0681257689044559
$ws.Cells.Item(7, 17).Value = 3.0618580948610665
$ws.Cells.Item(7, 19).Value = 3.0618581162154008
$ws.Cells.Item(7, 21).Value = 3.0618580948607583
$ws.Cells.Item(8, 5).Value = 0.95924193069612929
$ws.Cells.Item(8, 7).Value = 0.95924193115089085
$ws.Cells.Item(8, 9).Value = 0.95924193070224861
$ws.Cells.Item(8, 11).Value = 0.9572215172656795
$ws.Cells.Item(8, 13).Value = 0.957221517072259
$ws.Cells.Item(8, 15).Value = 0.957221517267422
$ws.Cells.Item(8, 17).Value = 0.96107338276594034
$ws.Cells.Item(8, 19).Value = 0.96107338456602553
$ws.Cells.Item(8, 21).Value = 0.96107338276642651
$ws.Cells.Item(9, 5).Value = 0.9490340568757426
$ws.Cells.Item(9, 7).Value = 0.94903405339453106
$ws.Cells.Item(9, 9).Value = 0.94903405688421294
$ws.Cells.Item(9, 11).Value = 0.94652627390687871
$ws.Cells.Item(9, 13).Value = 0.94652626888426039
$ws.Cells.Item(9, 15).Value = 0.94652627390882238
$ws.Cells.Item(9, 17).Value = 0.95130881612284057
$ws.Cells.Item(9, 19).Value = 0.95130881503115117
$ws.Cells.Item(9, 21).Value = 0.95130881612337159
$ws.Cells.Item(17, 5).Value = 2.1058245831317985
$ws.Cells.Item(18, 5).Value = 4.2894913950594482

# Add new row 19 (CWI / _ / 1), copying text style (quotePrefix) from row 18 labels
$ws.Range("A19").Value = "CWI"
$ws.Range("B19").Value = "_"
$ws.Range("C19").Value = 1
$ws.Range("A18:B18").Copy()
$ws.Range("A19:B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
